# #5: cash & deposit done
# Adds currency/cash/bank/deposit_type/deposit shared-string columns to the
# "現金" (cash) and "存款" (deposit) sheets, matching the updated exporter
# schema: owner,total,property_category,category,date,legislator_name,
# legislator_id,source_file,index appended after the existing columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "現金" (cash) -- 4th sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(4)

# Header row (row 1): B1/C1 keep the same header text, D1..K1 are brand new.
$ws.Cells.Item(1,2).Value = "currency"
$ws.Cells.Item(1,3).Value = "owner"

$ws.Cells.Item(1,3).Copy($ws.Cells.Item(1,4))
$ws.Cells.Item(1,4).Value = "total"
$ws.Cells.Item(1,3).Copy($ws.Cells.Item(1,5))
$ws.Cells.Item(1,5).Value = "property_category"
$ws.Cells.Item(1,3).Copy($ws.Cells.Item(1,6))
$ws.Cells.Item(1,6).Value = "category"
$ws.Cells.Item(1,3).Copy($ws.Cells.Item(1,7))
$ws.Cells.Item(1,7).Value = "date"
$ws.Cells.Item(1,3).Copy($ws.Cells.Item(1,8))
$ws.Cells.Item(1,8).Value = "legislator_name"
$ws.Cells.Item(1,3).Copy($ws.Cells.Item(1,9))
$ws.Cells.Item(1,9).Value = "legislator_id"
$ws.Cells.Item(1,3).Copy($ws.Cells.Item(1,10))
$ws.Cells.Item(1,10).Value = "source_file"
$ws.Cells.Item(1,3).Copy($ws.Cells.Item(1,11))
$ws.Cells.Item(1,11).Value = "index"

# Row 2 (index 49): D2 gets the amount that used to live in E2; E2 becomes
# the property_category literal "cash"; F2..K2 are new.
$ws.Cells.Item(2,4).Value = 2000000
$ws.Cells.Item(2,5).Value = "cash"
$ws.Cells.Item(2,4).Copy($ws.Cells.Item(2,6))
$ws.Cells.Item(2,6).Value = "normal"
$ws.Cells.Item(2,4).Copy($ws.Cells.Item(2,7))
$ws.Cells.Item(2,7).Value = "'2013-12-11"
$ws.Cells.Item(2,4).Copy($ws.Cells.Item(2,8))
$ws.Cells.Item(2,8).Value = "吳育仁"
$ws.Cells.Item(2,4).Copy($ws.Cells.Item(2,9))
$ws.Cells.Item(2,9).Value = 1734
$ws.Cells.Item(2,4).Copy($ws.Cells.Item(2,10))
$ws.Cells.Item(2,10).Value = "tmpbcc11"
$ws.Cells.Item(2,4).Copy($ws.Cells.Item(2,11))
$ws.Cells.Item(2,11).Value = 49

# Row 3 (index 50): D3 gets the value that used to live in E3; E3 becomes "cash".
$ws.Cells.Item(3,4).Value = 292000
$ws.Cells.Item(3,5).Value = "cash"
$ws.Cells.Item(3,4).Copy($ws.Cells.Item(3,6))
$ws.Cells.Item(3,6).Value = "normal"
$ws.Cells.Item(3,4).Copy($ws.Cells.Item(3,7))
$ws.Cells.Item(3,7).Value = "'2013-12-11"
$ws.Cells.Item(3,4).Copy($ws.Cells.Item(3,8))
$ws.Cells.Item(3,8).Value = "吳育仁"
$ws.Cells.Item(3,4).Copy($ws.Cells.Item(3,9))
$ws.Cells.Item(3,9).Value = 1734
$ws.Cells.Item(3,4).Copy($ws.Cells.Item(3,10))
$ws.Cells.Item(3,10).Value = "tmpbcc11"
$ws.Cells.Item(3,4).Copy($ws.Cells.Item(3,11))
$ws.Cells.Item(3,11).Value = 50

# Row 4 (index 51): D4 gets the value that used to live in E4; E4 becomes "cash".
$ws.Cells.Item(4,4).Value = 3000000
$ws.Cells.Item(4,5).Value = "cash"
$ws.Cells.Item(4,4).Copy($ws.Cells.Item(4,6))
$ws.Cells.Item(4,6).Value = "normal"
$ws.Cells.Item(4,4).Copy($ws.Cells.Item(4,7))
$ws.Cells.Item(4,7).Value = "'2013-12-11"
$ws.Cells.Item(4,4).Copy($ws.Cells.Item(4,8))
$ws.Cells.Item(4,8).Value = "吳育仁"
$ws.Cells.Item(4,4).Copy($ws.Cells.Item(4,9))
$ws.Cells.Item(4,9).Value = 1734
$ws.Cells.Item(4,4).Copy($ws.Cells.Item(4,10))
$ws.Cells.Item(4,10).Value = "tmpbcc11"
$ws.Cells.Item(4,4).Copy($ws.Cells.Item(4,11))
$ws.Cells.Item(4,11).Value = 51

# ---------------------------------------------------------------
# Sheet "存款" (deposit) -- 5th sheet
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(5)

# Header row (row 1): B1..E1 keep header text, F1..M1 are brand new.
$ws2.Cells.Item(1,2).Value = "bank"
$ws2.Cells.Item(1,3).Value = "deposit_type"
$ws2.Cells.Item(1,4).Value = "currency"
$ws2.Cells.Item(1,5).Value = "owner"

$ws2.Cells.Item(1,5).Copy($ws2.Cells.Item(1,6))
$ws2.Cells.Item(1,6).Value = "total"
$ws2.Cells.Item(1,5).Copy($ws2.Cells.Item(1,7))
$ws2.Cells.Item(1,7).Value = "property_category"
$ws2.Cells.Item(1,5).Copy($ws2.Cells.Item(1,8))
$ws2.Cells.Item(1,8).Value = "category"
$ws2.Cells.Item(1,5).Copy($ws2.Cells.Item(1,9))
$ws2.Cells.Item(1,9).Value = "date"
$ws2.Cells.Item(1,5).Copy($ws2.Cells.Item(1,10))
$ws2.Cells.Item(1,10).Value = "legislator_name"
$ws2.Cells.Item(1,5).Copy($ws2.Cells.Item(1,11))
$ws2.Cells.Item(1,11).Value = "legislator_id"
$ws2.Cells.Item(1,5).Copy($ws2.Cells.Item(1,12))
$ws2.Cells.Item(1,12).Value = "source_file"
$ws2.Cells.Item(1,5).Copy($ws2.Cells.Item(1,13))
$ws2.Cells.Item(1,13).Value = "index"

function Fill-DepositRow($row, $total, $owner, $index) {
    # F<row> gets the amount that used to be in G<row>; G<row> becomes "deposit".
    $ws2.Cells.Item($row,6).Value = $total
    $ws2.Cells.Item($row,7).Value = "deposit"
    $ws2.Cells.Item($row,6).Copy($ws2.Cells.Item($row,8))
    $ws2.Cells.Item($row,8).Value = "normal"
    $ws2.Cells.Item($row,6).Copy($ws2.Cells.Item($row,9))
    $ws2.Cells.Item($row,9).Value = "'2013-12-11"
    $ws2.Cells.Item($row,6).Copy($ws2.Cells.Item($row,10))
    $ws2.Cells.Item($row,10).Value = "吳育仁"
    $ws2.Cells.Item($row,6).Copy($ws2.Cells.Item($row,11))
    $ws2.Cells.Item($row,11).Value = 1734
    $ws2.Cells.Item($row,6).Copy($ws2.Cells.Item($row,12))
    $ws2.Cells.Item($row,12).Value = "tmpbcc11"
    $ws2.Cells.Item($row,6).Copy($ws2.Cells.Item($row,13))
    $ws2.Cells.Item($row,13).Value = $index
}

Fill-DepositRow 2 1023635 "吳育仁" 56
Fill-DepositRow 3 103269 "吳育仁" 57
Fill-DepositRow 4 404136 "吳育仁" 58
Fill-DepositRow 5 1179285 "蔡瓊姿" 59
Fill-DepositRow 6 133337 "吳育仁" 60
Fill-DepositRow 7 1070331 "吳育仁" 61
Fill-DepositRow 8 1974 "吳育仁" 62
Fill-DepositRow 9 866838 "蔡瓊姿" 64
Fill-DepositRow 10 16050.48 "蔡瓊姿" 65
Fill-DepositRow 11 2923 "蔡瓊姿" 66
